# "Added student dashboard on teacher's panel"
#
# Applies the Tasks.xlsx edits: a couple of new "completed" percentage
# markers, a task-owner name fix, three brand-new task/issue rows, the
# removal of a few obsolete rows (finMind/TASK/Fix wizard term-perm/Lain)
# and the re-shuffling that leaves behind, plus the resulting selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Mark two existing tasks as 100% complete (new G3 / G13 cells) -----
# ---------------------------------------------------------------------
$ws.Range("G3").Value = 1
$ws.Range("G3").NumberFormat = $ws.Range("G11").NumberFormat
$ws.Range("G3").Font.Color = $ws.Range("G11").Font.Color

$ws.Range("G13").Value = 1
$ws.Range("G13").NumberFormat = $ws.Range("G11").NumberFormat
$ws.Range("G13").Font.Color = $ws.Range("G11").Font.Color

# ---------------------------------------------------------------------
# 2) Row 12's "Fix By" owner was corrected from "Lain" to "lian" -------
# ---------------------------------------------------------------------
$ws.Range("C12").Value = "lian"

# ---------------------------------------------------------------------
# 3) Replace the old blank rows 20-21 and the "finMind" row 22 with ----
#    the new student-dashboard task plus two bug/issue rows ------------
# ---------------------------------------------------------------------
$ws.Range("A20").Value = "Task  10"
$ws.Range("B20").Value = "Add student name tag on top of student cosole "
$ws.Range("C20").Value = "Lian"
$ws.Range("E20").Value = "Add a name tag to identify student. Especially useful when multiple students are present"
$ws.Range("F20").Value = (Get-Date -Year 2022 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F20").NumberFormat = "d-mmm"
$ws.Rows.Item(20).RowHeight = 45

$ws.Range("A21").Value = "Issue 1"
$ws.Range("B21").Value = "User Name should be w/o special chars since we use it as video id"
$ws.Range("C21").Value = "Lian"
$ws.Rows.Item(21).RowHeight = 30

$ws.Range("A22").Value = "Issue 2"
$ws.Range("B22").Value = "Timer not working every time"
$ws.Range("C22").Value = "Lian"
$ws.Rows.Item(22).RowHeight = 15

# ---------------------------------------------------------------------
# 4) The old "TASK / ADD ERROR PAGE" row (23), "Fix wizard term-perm ---
#    option" row (24) and "From application to pipeline" row (25) -----
#    are cleared out, leaving just the section-header style in A23:A25 -
# ---------------------------------------------------------------------
$ws.Range("A23").ClearContents()
$ws.Range("B23").ClearContents()
$ws.Range("A23").Font.Color = $ws.Range("A22").Font.Color
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").Font.Italic = $true
$ws.Range("A23").Font.Size = 14
$ws.Rows.Item(23).RowHeight = 18.75

$ws.Range("A24").ClearContents()
$ws.Range("B24").ClearContents()
$ws.Range("A24").Font.Color = $ws.Range("A22").Font.Color
$ws.Range("A24").Font.Bold = $true
$ws.Range("A24").Font.Italic = $true
$ws.Range("A24").Font.Size = 14
$ws.Rows.Item(24).RowHeight = 18.75

$ws.Range("A25").ClearContents()
$ws.Range("B25").ClearContents()
$ws.Range("A25").Font.Color = $ws.Range("A22").Font.Color
$ws.Range("A25").Font.Bold = $true
$ws.Range("A25").Font.Italic = $true
$ws.Range("A25").Font.Size = 14
$ws.Rows.Item(25).RowHeight = 18.75

# ---------------------------------------------------------------------
# 5) The two surviving detail lines move down to rows 26 and 28 --------
#    (row 27 stays empty), now living in column B instead of A/B ------
# ---------------------------------------------------------------------
$ws.Range("B26").Value = "ADD ERROR PAGE (MESSAGES)"
$ws.Range("B26").WrapText = $true

$ws.Range("B28").Value = "From application to pipeline"
$ws.Range("B28").WrapText = $true

# ---------------------------------------------------------------------
# 6) Final selection left behind in the sheet --------------------------
# ---------------------------------------------------------------------
$ws.Range("E12").Select()
